# Updated omikron excel, now using numbers from daily report instead
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: Dato, AntalTest, AntalOmikron for each date row (rows 2-18)
$data = @(
    @("2021-11-22", 4470, 1),
    @("2021-11-23", 4666, 1),
    @("2021-11-24", 3982, 1),
    @("2021-11-25", 4059, 4),
    @("2021-11-26", 4114, 7),
    @("2021-11-27", 3813, 3),
    @("2021-11-28", 3849, 10),
    @("2021-11-29", 5048, 11),
    @("2021-11-30", 5368, 25),
    @("2021-12-01", 4491, 76),
    @("2021-12-02", 4528, 60),
    @("2021-12-03", 5126, 77),
    @("2021-12-04", 5058, 101),
    @("2021-12-05", 4768, 170),
    @("2021-12-06", 7028, 356),
    @("2021-12-07", 7162, 581),
    @("2021-12-08", 1759, 311)
)

# Header row stays the same, but re-write to be safe
$ws.Cells.Item(1, 1).Value = "Dato"
$ws.Cells.Item(1, 2).Value = "AntalTest"
$ws.Cells.Item(1, 3).Value = "AntalOmikron"
$ws.Cells.Item(1, 4).Value = "Ratio"

$row = 2
foreach ($entry in $data) {
    $dateStr = $entry[0]
    $antalTest = $entry[1]
    $antalOmikron = $entry[2]

    # Leading apostrophe forces text entry (matches original quotePrefix style)
    $ws.Cells.Item($row, 1).Value = "'" + $dateStr
    $ws.Cells.Item($row, 2).Value = $antalTest
    $ws.Cells.Item($row, 3).Value = $antalOmikron
    $ws.Cells.Item($row, 4).Formula = "=100*C$row/B$row"

    $row = $row + 1
}

# Remaining rows 19-23 just keep the existing (date) formatting on column A, empty
for ($r = 19; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
}

# Rows 24-25 did not exist before; copy the formatting from row 19 (style only)
# down into them so the sheet's used range grows to A1:D25, matching target.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A24:A25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F5").Select()
